# Kitabevi Modellemesi.xlsx - "Add files via upload" edit
# Adds a second worksheet ("Sayfa2") with three explanatory sections
# (SÖZLÜK / REHBER, ANA TABLOLAR, HAREKETLİLİK) and makes it the active tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- create Sayfa2 right after Sayfa1 ------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sayfa2"

# --- write the text content in the same order the original author did ----
# (controls shared-string allocation order: 69..77)
$ws2.Range("B3").Value  = 'İl adı, ilçe adı, eğitim durumu, okul, üniversite adı, çalışılan pozisyon…..'
$ws2.Range("A6").Value  = '2. ANA TABLOLAR /  GENEL TABLOLAR'
$ws2.Range("B8").Value  = 'Kitaplar, Ürünler, Personel, Müşteriler, Tedarikçiler, Hesaplar…..'
$ws2.Range("A11").Value = '3. HAREKETLİLİK TABLOLARI'
$ws2.Range("A1").Value  = '1. SÖZLÜK / REHBER TABLOLARI'
$ws2.Range("B13").Value = 'Girişler, Satışlar, Not Girişleri, ToplantıTakibi, Gönderilen Mesajlar,…..'
$ws2.Range("A2").Value  = 'İçinde çoğunlukla id no ve o id noya ait bir içerik barındıran, çok nadir 3 .alana sahip olan tablolardır. Yaratılma sebebi , başka tablolarda terkrarlanan alan verilerini hatalı girişlere ve veri kayıplarına karşın bir tabloda yedeklenmesini sağlamaktır. Bu tabloya gerektiğinde yabancıl anahtar alan ( foreign key) ile referans için başvurulur.  '
$ws2.Range("A7").Value  = 'Bir ürün, kişi, kurum, bina gibi kavramlara ait , o kavramı oluşturan ve o kavramı tamamlayacak tüm bilgi alanlarının verilerinin tutulduğu tablolardır. Bu tür tablolar çoğunlukla veritabanı yapısının / modellemesinin ana enstrumanıdır. Veri girişi yapılırken eğer içlerindeki kayıt satırlarında  tekrarlanma ihtimali yüksek  veriler olacaksa bu veriler ilgili sözlük tablolarından foreign key kısıtlaması getirilerek çekilecek referans veriler ile temsil edilmelidir. Bu durum VERİLERİN TUTARLILIĞI ile isimlendirilir. Aksi duruma ise VERİLER ARASI TUTARSIZLIK ( inconsistency )  denir.'
$ws2.Range("A12").Value = 'Veritabanının oluşturulmasına amaç olan hareketlerin takibi için yaratılan tablolardır. Bu tür tablolara kayıt işlemi çoğunlukla hemen hemen her gün defalarca yapılır. Kayıt işlemi sırasında çoğu alanlarının veri girişi foreign key kısıtlaması ile referans alınan diğer tablolar ile ilişkilendirilerek yapılır. Bu tablolar genellikle ana tablolar olur. Veritabanının kapladığı alan içinde en büyük payı bu tür tablolar alır. Çoğunlukla üzerlerinde en fazla sorgu işleminin yapıldığı tablolardır.'

# --- section headers: bold, default colour (A1 / A6 / A11) --------------
$ws2.Range("A1").Font.Bold = $true
$ws2.Range("A6").Font.Bold = $true
$ws2.Range("A11").Font.Bold = $true

# --- sub-bullets: blue text (B3/B4, B8/B9, B13/B14) ----------------------
$ws2.Range("B3:B4").Font.Color = 12611584
$ws2.Range("B8:B9").Font.Color = 12611584
$ws2.Range("B13:B14").Font.Color = 12611584

# --- long description paragraphs: merged, wrapped, vertically centred ---
$descRanges = @("A2:Q2", "A7:Q7", "A12:Q12")
foreach ($r in $descRanges) {
    $rng = $ws2.Range($r)
    $rng.Merge()
    $rng.HorizontalAlignment = -4131
    $rng.VerticalAlignment = -4108
    $rng.WrapText = $true
}
$ws2.Rows.Item(2).RowHeight = 43.8
$ws2.Rows.Item(7).RowHeight = 68.4
$ws2.Rows.Item(12).RowHeight = 52.8

# --- separator rows: "Good"/İyi built-in style (green fill/text) --------
$ws2.Range("A5:Q5").Style = "Good"
$ws2.Range("A10:Q10").Style = "Good"
$ws2.Range("A15:Q15").Style = "Good"

# --- sheet view ------------------------------------------------------------
$ws2.Activate()
$excel.ActiveWindow.Zoom = 120
$ws2.Range("A13").Select()

Write-Host "Sayfa2 added with" $wb.Worksheets.Count "worksheets total"
